# daily auto push: 2026-01-31 18:49 UTC
# Insert a new data row for 2026/01/31 just before the existing row 748
# (2026/12/29), shifting rows 748:789 down to 749:790.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 748, pushing row 748 (and everything
# below it) down by one row.
$ws.Rows("748:748").Insert()

# Populate the newly inserted row with the new record. Force the date
# column to plain text first so Excel does not auto-convert the
# "yyyy/mm/dd" string into a date serial value, then restore the cell's
# default styling so it matches its neighbours (no explicit style).
$ws.Cells.Item(748, 1).NumberFormat = "@"
$ws.Cells.Item(748, 1).Value = "2026/01/31"
$ws.Cells.Item(748, 1).Style = "Normal"
$ws.Cells.Item(748, 2).Value = "土"
$ws.Cells.Item(748, 3).Value = 23
$ws.Cells.Item(748, 4).Value = 201
